# Updated cryptos list values per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a Excel number
# (losing the exact decimal text, e.g. trailing zeros / scientific notation).
# These are written with a leading apostrophe so Excel keeps them as text,
# exactly like the original workbook already stores D2/D3/D7/... etc.
$forceText = @(
    "D5", "D6", "D8", "D9", "D11", "D12", "D13", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D31", "D32", "D34", "D36", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)

$updates = @{
    "D2" = "67.341.48"
    "E2" = "  -0.09%  "
    "D3" = "3.498.71"
    "E3" = "  -0.21%  "
    "E4" = "  +0.04%  "
    "D5" = "606.85"
    "E5" = "  -0.13%  "
    "D6" = "151.35"
    "E6" = "  +1.23%  "
    "D7" = "3.499.59"
    "E7" = "  -0.16%  "
    "D8" = "1.00"
    "E8" = "  -0.02%  "
    "D9" = "0.486"
    "E9" = "  +1.00%  "
    "E10" = "  +2.76%  "
    "D11" = "7.48"
    "E11" = "  +6.07%  "
    "D12" = "0.431"
    "E12" = "  +1.71%  "
    "D13" = "32.45"
    "E13" = "  +2.42%  "
    "E14" = "  -1.89%  "
    "D15" = "4.085.60"
    "E15" = "  -0.29%  "
    "B16" = "WrappedEther"
    "C16" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D16" = "3.498.21"
    "E16" = "  -0.19%  "
    "B17" = "WrappedBTC"
    "C17" = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
    "D17" = "67.307.13"
    "E17" = "  -0.11%  "
    "E18" = "  -0.55%  "
    "D19" = "6.53"
    "E19" = "  +1.91%  "
    "D20" = "15.48"
    "E20" = "  +2.14%  "
    "D21" = "9.81"
    "E21" = "  +6.42%  "
    "D22" = "446.17"
    "E22" = "  +0.02%  "
    "D23" = "0.631"
    "E23" = "  +1.32%  "
    "D24" = "77.80"
    "E24" = "  +0.45%  "
    "D25" = "3.634.03"
    "E25" = "  -0.29%  "
    "B26" = "PEPE"
    "C26" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D26" = "0.0000127"
    "E26" = "  -0.34%  "
    "B27" = "Dai"
    "C27" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D27" = "1.00"
    "E27" = "  +0.04%  "
    "D28" = "8.83"
    "E28" = "  +6.04%  "
    "D29" = "10.05"
    "E29" = "  -2.13%  "
    "E30" = "  +0.54%  "
    "D31" = "1.64"
    "E31" = "  +6.39%  "
    "D32" = "0.169"
    "E32" = "  +2.91%  "
    "E33" = "  -0.15%  "
    "D34" = "25.68"
    "E34" = "  +0.09%  "
    "E35" = "  +0.98%  "
    "D36" = "1.87"
    "E36" = "  +1.86%  "
    "D37" = "3.488.70"
    "E37" = "  -0.23%  "
    "E38" = "  -0.27%  "
    "E39" = "  +0.00%  "
    "E40" = "  +5.78%  "
    "D41" = "0.999"
    "E41" = "  +0.00%  "
    "D42" = "174.41"
    "E42" = "  -1.21%  "
    "D43" = "0.0895"
    "E43" = "  +2.79%  "
    "D44" = "5.46"
    "E44" = "  +0.80%  "
    "D45" = "30.09"
    "E45" = "  +10.97%  "
    "D46" = "0.876"
    "E46" = "  -0.35%  "
    "D47" = "47.05"
    "E47" = "  +3.68%  "
    "D48" = "1.30"
    "E48" = "  +3.07%  "
    "B49" = "Cosmos"
    "C49" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "D49" = "7.65"
    "E49" = "  +1.14%  "
    "B50" = "dogwifhat"
    "C50" = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
    "D50" = "2.51"
    "E50" = "  -2.30%  "
    "D51" = "0.253"
    "E51" = "  +3.19%  "
}

foreach ($cell in $updates.Keys) {
    $value = $updates[$cell]
    if ($forceText -contains $cell) {
        $value = "'" + $value
    }
    $ws.Range($cell).Value = $value
}
